$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for rule R10 (cell E8)
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active cell selection on the sheet
$ws.Range("E8").Select()
